$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 48
$ws.Range("H48").Value = 10000
$ws.Range("I48").Value = 10000
$ws.Range("K48").Value = 30000
$ws.Range("M48").Value = -29708
# Row 56
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 30000
$ws.Range("M56").Value = -29466
# Row 106
$ws.Range("H106").Value = 5130477
$ws.Range("I106").Value = 5557933.5
$ws.Range("K106").Value = 5557933.5
$ws.Range("M106").Value = -5557302.5
# Row 129
$ws.Range("H129").Value = 1742.6428
$ws.Range("I129").Value = 1429.4
$ws.Range("K129").Value = 4288.200000000001
$ws.Range("M129").Value = 711.7999999999993
# Row 131
$ws.Range("H131").Value = 11003920
$ws.Range("I131").Value = 25002300
$ws.Range("K131").Value = 75006900
$ws.Range("M131").Value = -75001860
# Row 132
$ws.Range("H132").Value = 1807.2354
$ws.Range("I132").Value = 1857.5
$ws.Range("K132").Value = 5572.5
$ws.Range("M132").Value = -3042.5
# Row 137
$ws.Range("H137").Value = 29155.367
$ws.Range("I137").Value = 42958.05
$ws.Range("K137").Value = 128874.15
$ws.Range("M137").Value = -126324.15

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2166.9565
$ws.Range("I2").Value = 2189.1052
$ws.Range("K2").Value = 2189.1052
$ws.Range("M2").Value = -2076.1052
# Row 32
$ws.Range("H32").Value = 29415520
$ws.Range("I32").Value = 40002908
$ws.Range("K32").Value = 40002908
$ws.Range("M32").Value = -40002621
# Row 61
$ws.Range("H61").Value = 3468.0454
$ws.Range("I61").Value = 3449.75
$ws.Range("K61").Value = 3449.75
$ws.Range("M61").Value = -3237.75
# Row 63
$ws.Range("H63").Value = 4624.875
$ws.Range("I63").Value = 2500
$ws.Range("J63").Value = 5333.1665
$ws.Range("K63").Value = 2500
$ws.Range("L63").Value = 5333.1665
$ws.Range("M63").Value = -1814
$ws.Range("N63").Value = -6705.1665
# Row 66
$ws.Range("H66").Value = 4624.875
$ws.Range("I66").Value = 2500
$ws.Range("J66").Value = 5333.1665
$ws.Range("K66").Value = 12500
$ws.Range("L66").Value = 26665.8325
$ws.Range("M66").Value = -9068
$ws.Range("N66").Value = -33529.8325
# Row 74
$ws.Range("H74").Value = 2823.1667
$ws.Range("I74").Value = 2402.6177
$ws.Range("J74").Value = 4610.5
$ws.Range("K74").Value = 2402.6177
$ws.Range("L74").Value = 4610.5
$ws.Range("M74").Value = -1528.6177
$ws.Range("N74").Value = -6358.5
# Row 77
$ws.Range("H77").Value = 2823.1667
$ws.Range("I77").Value = 2402.6177
$ws.Range("J77").Value = 4610.5
$ws.Range("K77").Value = 12013.0885
$ws.Range("L77").Value = 23052.5
$ws.Range("M77").Value = -7645.088499999998
$ws.Range("N77").Value = -31788.5
# Row 110
$ws.Range("H110").Value = 3497.25
$ws.Range("I110").Value = 3496.3333
$ws.Range("K110").Value = 3496.3333
$ws.Range("M110").Value = -1451.3333
# Row 116
$ws.Range("H116").Value = 2166.9565
$ws.Range("I116").Value = 2189.1052
$ws.Range("K116").Value = 2189.1052
$ws.Range("M116").Value = 104.8948
# Row 132
$ws.Range("H132").Value = 3276.3142
$ws.Range("I132").Value = 2931.8147
$ws.Range("K132").Value = 8795.444100000001
$ws.Range("M132").Value = -6265.444100000001
# Row 136
$ws.Range("H136").Value = 3468.0454
$ws.Range("I136").Value = 3449.75
$ws.Range("K136").Value = 10349.25
$ws.Range("M136").Value = -7799.25

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2166.9565
$ws.Range("I3").Value = 2189.1052
$ws.Range("K3").Value = 2189.1052
$ws.Range("M3").Value = -2075.1052
# Row 7
$ws.Range("H7").Value = 4003420
$ws.Range("I7").Value = 5001775
$ws.Range("K7").Value = 5001775
$ws.Range("M7").Value = -5001662
# Row 107
$ws.Range("H107").Value = 3478.375
$ws.Range("I107").Value = 3974.5
$ws.Range("J107").Value = 3313
$ws.Range("K107").Value = 3974.5
$ws.Range("L107").Value = 3313
$ws.Range("M107").Value = -2054.5
$ws.Range("N107").Value = -7153

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2575.4644
$ws.Range("I58").Value = 2401.625
$ws.Range("J58").Value = 3618.5
$ws.Range("K58").Value = 2401.625
$ws.Range("L58").Value = 3618.5
$ws.Range("M58").Value = -2198.625
$ws.Range("N58").Value = -4024.5
# Row 132
$ws.Range("H132").Value = 2842.9211
$ws.Range("I132").Value = 2632.3428
$ws.Range("K132").Value = 7897.028399999999
$ws.Range("M132").Value = -5367.028399999999
# Row 136
$ws.Range("H136").Value = 2575.4644
$ws.Range("I136").Value = 2401.625
$ws.Range("J136").Value = 3618.5
$ws.Range("K136").Value = 7204.875
$ws.Range("L136").Value = 10855.5
$ws.Range("M136").Value = -4654.875
$ws.Range("N136").Value = -15955.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 70
$ws.Range("H70").Value = 4994.5
$ws.Range("I70").Value = 4994.5
$ws.Range("K70").Value = 14983.5
$ws.Range("M70").Value = -14668.5
# Row 73
$ws.Range("H73").Value = 4994.5
$ws.Range("I73").Value = 4994.5
$ws.Range("K73").Value = 14983.5
$ws.Range("M73").Value = -13891.5
# Row 75
$ws.Range("H75").Value = 24886.5
$ws.Range("I75").Value = 5273.5
$ws.Range("J75").Value = 44499.5
$ws.Range("K75").Value = 15820.5
$ws.Range("L75").Value = 133498.5
$ws.Range("M75").Value = -14822.5
$ws.Range("N75").Value = -135494.5
# Row 78
$ws.Range("H78").Value = 24886.5
$ws.Range("I78").Value = 5273.5
$ws.Range("J78").Value = 44499.5
$ws.Range("K78").Value = 47461.5
$ws.Range("L78").Value = 400495.5
$ws.Range("M78").Value = -42469.5
$ws.Range("N78").Value = -410479.5
# Row 121
$ws.Range("H121").Value = 5039656
$ws.Range("I121").Value = 544.6667
$ws.Range("K121").Value = 1634.0001
$ws.Range("M121").Value = -324.0001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5999.9614
$ws.Range("I70").Value = 11374.75
$ws.Range("J70").Value = 5022.727
$ws.Range("K70").Value = 11374.75
$ws.Range("L70").Value = 5022.727
$ws.Range("M70").Value = -11104.75
$ws.Range("N70").Value = -5562.727
# Row 73
$ws.Range("H73").Value = 5999.9614
$ws.Range("I73").Value = 11374.75
$ws.Range("J73").Value = 5022.727
$ws.Range("K73").Value = 11374.75
$ws.Range("L73").Value = 5022.727
$ws.Range("M73").Value = -10438.75
$ws.Range("N73").Value = -6894.727
# Row 132
$ws.Range("H132").Value = 6767.8335
$ws.Range("I132").Value = 6701.5713
$ws.Range("K132").Value = 20104.7139
$ws.Range("M132").Value = -17574.7139

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 18522340
$ws.Range("J40").Value = 7997
$ws.Range("L40").Value = 7997
$ws.Range("N40").Value = -8269

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5194.3076
$ws.Range("I62").Value = 3555
$ws.Range("J62").Value = 5922.8887
$ws.Range("K62").Value = 3555
$ws.Range("L62").Value = 5922.8887
$ws.Range("M62").Value = -2931
$ws.Range("N62").Value = -7170.8887
# Row 65
$ws.Range("H65").Value = 5194.3076
$ws.Range("I65").Value = 3555
$ws.Range("J65").Value = 5922.8887
$ws.Range("K65").Value = 17775
$ws.Range("L65").Value = 29614.4435
$ws.Range("M65").Value = -14655
$ws.Range("N65").Value = -35854.4435
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()
# Row 95
$ws.Range("H95").Value = 78663.336
$ws.Range("J95").Value = 78663.336
$ws.Range("L95").Value = 78663.336
$ws.Range("N95").Value = -84155.336
# Row 113
$ws.Range("H113").Value = 1519.2
$ws.Range("I113").Value = 1483.6666
$ws.Range("J113").Value = 1572.5
$ws.Range("K113").Value = 4450.9998
$ws.Range("L113").Value = 4717.5
$ws.Range("M113").Value = -2280.9998
$ws.Range("N113").Value = -9057.5
# Row 135
$ws.Range("H135").Value = 49995
$ws.Range("J135").Value = 49995
$ws.Range("L135").Value = 49995
$ws.Range("N135").Value = -60135

